$d = $word.ActiveDocument

# --- Edit 1 ---------------------------------------------------------------
# Merge the leading "            " (12-space) whitespace-only run into the
# following "private static final " run (they share identical rPr), so the
# two runs become a single run with the combined text. Only the first
# occurrence (the FRAME_WIDTH declaration, preceded by the 12-space run)
# should be touched -- the second "private static final" (FRAME_HEIGHT,
# preceded by a tab) must stay untouched. Scope the Find to the specific
# paragraph so only that occurrence is affected.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "            private static final int FRAME_WIDTH = 400;") {
        $r = $p.Range
        $r.Find.Execute("            private static final ", $true, $false, $false, $false, $false, $true, 1, $false, "            private static final ", 2) | Out-Null
        break
    }
}

# --- Edit 2 ---------------------------------------------------------------
# Remove the "Workout Selection Class:" section entirely: the page break
# paragraph, the "Workout Selection Class:" paragraph, the blank paragraph,
# and the paragraph holding the WorkOutSelection class-diagram textbox plus
# the "Class Attribute List:" run. Keep the final (now-empty) paragraph that
# carries the _GoBack bookmark, but strip its paragraph formatting too.
$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd()
    if ($startIdx -eq -1 -and $t -eq "Workout Selection Class:") {
        $startIdx = $i - 1
    }
    if ($startIdx -ne -1 -and $i -gt $startIdx -and $t -eq "Class Attribute List:") {
        $endIdx = $i
        break
    }
}

if ($startIdx -ne -1 -and $endIdx -ne -1) {
    $startPara = $d.Paragraphs($startIdx)
    $endPara = $d.Paragraphs($endIdx)
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()

    # After the delete, the paragraph that used to be at $startIdx (the page
    # break paragraph) is now the merged/last paragraph in the document --
    # it carries the surviving bookmark. Strip its paragraph formatting.
    $lastPara = $d.Paragraphs($startIdx)
    $lastPara.Range.Select()
    $word.Selection.ClearFormatting()
}
